$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-20 -> 2023-09-21, i.e. Excel serial 45189 -> 45190) for every
# data row (rows 2 through 224).
$ws.Range("C2:C224").Value = 45190
